$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.414.44"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "1.853.68"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.74"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4681"
$ws.Range("E7").Value = "  -1.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2733"
$ws.Range("E8").Value = "  -0.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06310"
$ws.Range("E9").Value = "  -2.11%  "

$ws.Range("D10").Value = "1.852.25"
$ws.Range("E10").Value = "  -1.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07460"
$ws.Range("E11").Value = "  +0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.64"
$ws.Range("E12").Value = "  +3.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.949"
$ws.Range("E13").Value = "  -1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.10"
$ws.Range("E14").Value = "  -1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6212"
$ws.Range("E15").Value = "  -2.20%  "

$ws.Range("D16").Value = "30.353.64"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.31"
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007302"
$ws.Range("E19").Value = "  -0.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.42"
$ws.Range("E20").Value = "  -3.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9988"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.931"
$ws.Range("E22").Value = "  -3.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.893"
$ws.Range("E23").Value = "  -2.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "167.17"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.178"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.88"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.878"
$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1022"
$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.373"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.095"
$ws.Range("E30").Value = "  -3.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.819"
$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04896"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.140"
$ws.Range("E33").Value = "  -1.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7049"
$ws.Range("E34").Value = "  -3.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.698"
$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01935"
$ws.Range("E36").Value = "  -1.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.666"
$ws.Range("E37").Value = "  +1.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8740"
$ws.Range("E38").Value = "  -3.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.952"
$ws.Range("E39").Value = "  -2.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.59"
$ws.Range("E40").Value = "  -0.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9996"
$ws.Range("E41").Value = "  +0.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.521"
$ws.Range("E42").Value = "  -1.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4055"
$ws.Range("E43").Value = "  -1.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.075"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.28"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1211"
$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.589"
$ws.Range("E47").Value = "  -2.17%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.42"
$ws.Range("E48").Value = "  +1.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05547"
$ws.Range("E49").Value = "  -0.80%  "

$ws.Range("E50").Value = "  -3.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3667"
$ws.Range("E51").Value = "  -1.48%  "
